$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 9103334
$ws.Range("I33").Value = 24330.2
$ws.Range("J33").Value = 18182338
$ws.Range("K33").Value = 24330.2
$ws.Range("L33").Value = 18182338
$ws.Range("M33").Value = -24101.2
$ws.Range("N33").Value = -18182796
$ws.Range("H40").Value = 6758696
$ws.Range("I40").Value = 1960.9062
$ws.Range("K40").Value = 1960.9062
$ws.Range("M40").Value = -1785.9062
$ws.Range("H137").Value = 1704.3462
$ws.Range("I137").Value = 1491.5
$ws.Range("J137").Value = 2875
$ws.Range("K137").Value = 4474.5
$ws.Range("L137").Value = 8625
$ws.Range("M137").Value = -1924.5
$ws.Range("N137").Value = -13725

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1037.7
$ws.Range("I2").Value = 670.5
$ws.Range("J2").Value = 2506.5
$ws.Range("K2").Value = 670.5
$ws.Range("L2").Value = 2506.5
$ws.Range("M2").Value = -557.5
$ws.Range("N2").Value = -2732.5
$ws.Range("H49").Value = 13040
$ws.Range("J49").Value = 13040
$ws.Range("L49").Value = 13040
$ws.Range("N49").Value = -13560
$ws.Range("H74").Value = 1715.1818
$ws.Range("I74").Value = 1655
$ws.Range("J74").Value = 1759.5264
$ws.Range("K74").Value = 1655
$ws.Range("L74").Value = 1759.5264
$ws.Range("M74").Value = -781
$ws.Range("N74").Value = -3507.5264
$ws.Range("H77").Value = 1715.1818
$ws.Range("I77").Value = 1655
$ws.Range("J77").Value = 1759.5264
$ws.Range("K77").Value = 8275
$ws.Range("L77").Value = 8797.632
$ws.Range("M77").Value = -3907
$ws.Range("N77").Value = -17533.632
$ws.Range("H116").Value = 1037.7
$ws.Range("I116").Value = 670.5
$ws.Range("J116").Value = 2506.5
$ws.Range("K116").Value = 670.5
$ws.Range("L116").Value = 2506.5
$ws.Range("M116").Value = 1623.5
$ws.Range("N116").Value = -7094.5

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1037.7
$ws.Range("I3").Value = 670.5
$ws.Range("J3").Value = 2506.5
$ws.Range("K3").Value = 670.5
$ws.Range("L3").Value = 2506.5
$ws.Range("M3").Value = -556.5
$ws.Range("N3").Value = -2734.5
$ws.Range("H20").Value = 50014480
$ws.Range("I20").Value = 100001750
$ws.Range("K20").Value = 100001750
$ws.Range("M20").Value = -100001503
$ws.Range("H107").Value = 1176.84
$ws.Range("I107").Value = 1083.6666
$ws.Range("J107").Value = 1262.8462
$ws.Range("K107").Value = 1083.6666
$ws.Range("L107").Value = 1262.8462
$ws.Range("M107").Value = 836.3334
$ws.Range("N107").Value = -5102.8462
$ws.Range("H134").Value = 4742.8945
$ws.Range("I134").Value = 5574.3076
$ws.Range("K134").Value = 16722.9228
$ws.Range("M134").Value = -14187.9228

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5489.92
$ws.Range("I31").Value = 1883.6552
$ws.Range("J31").Value = 10470
$ws.Range("K31").Value = 1883.6552
$ws.Range("L31").Value = 10470
$ws.Range("M31").Value = -1588.6552
$ws.Range("N31").Value = -11060
$ws.Range("H34").Value = 5489.92
$ws.Range("I34").Value = 1883.6552
$ws.Range("J34").Value = 10470
$ws.Range("K34").Value = 1883.6552
$ws.Range("L34").Value = 10470
$ws.Range("M34").Value = -1681.6552
$ws.Range("N34").Value = -10874
$ws.Range("H62").Value = 5098.25
$ws.Range("I62").Value = 5983.2856
$ws.Range("J62").Value = 3859.2
$ws.Range("K62").Value = 5983.2856
$ws.Range("L62").Value = 3859.2
$ws.Range("M62").Value = -5359.2856
$ws.Range("N62").Value = -5107.2
$ws.Range("H65").Value = 5098.25
$ws.Range("I65").Value = 5983.2856
$ws.Range("J65").Value = 3859.2
$ws.Range("K65").Value = 29916.428
$ws.Range("L65").Value = 19296
$ws.Range("M65").Value = -26796.428
$ws.Range("N65").Value = -25536
$ws.Range("H68").Value = 28813.5
$ws.Range("J68").Value = 28813.5
$ws.Range("L68").Value = 28813.5
$ws.Range("N68").Value = -30311.5
$ws.Range("H71").Value = 28813.5
$ws.Range("J71").Value = 28813.5
$ws.Range("L71").Value = 86440.5
$ws.Range("N71").Value = -93928.5
$ws.Range("I86").Value = 250002480
$ws.Range("J86").Value = 2460
$ws.Range("K86").Value = 250002480
$ws.Range("L86").Value = 2460
$ws.Range("M86").Value = -250001357
$ws.Range("N86").Value = -4706
$ws.Range("I89").Value = 250002480
$ws.Range("J89").Value = 2460
$ws.Range("K89").Value = 1250012400
$ws.Range("L89").Value = 12300
$ws.Range("M89").Value = -1250006784
$ws.Range("N89").Value = -23532

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H49").Value = 15200
$ws.Range("J49").Value = 15200
$ws.Range("L49").Value = 15200
$ws.Range("N49").Value = -15568
$ws.Range("H126").Value = 5868.3213
$ws.Range("I126").Value = 7674.706
$ws.Range("K126").Value = 23024.118
$ws.Range("M126").Value = -20554.118

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1229.9333
$ws.Range("J46").Value = 1274.9
$ws.Range("L46").Value = 1274.9
$ws.Range("N46").Value = -1650.9
$ws.Range("H68").Value = 111113064
$ws.Range("I68").Value = 2185.7144
$ws.Range("J68").Value = 500001150
$ws.Range("K68").Value = 2185.7144
$ws.Range("L68").Value = 500001150
$ws.Range("M68").Value = -1436.7144
$ws.Range("N68").Value = -500002648
$ws.Range("H71").Value = 111113064
$ws.Range("I71").Value = 2185.7144
$ws.Range("J71").Value = 500001150
$ws.Range("K71").Value = 10928.572
$ws.Range("L71").Value = 2500005750
$ws.Range("M71").Value = -7184.572
$ws.Range("N71").Value = -2500013238
$ws.Range("H133").Value = 73725.2
$ws.Range("J133").Value = 73725.2
$ws.Range("L133").Value = 73725.2
$ws.Range("N133").Value = -78785.2
$ws.Range("H136").Value = 7600.0977
$ws.Range("I136").Value = 5308.8823
$ws.Range("J136").Value = 18728.857
$ws.Range("K136").Value = 15926.6469
$ws.Range("L136").Value = 56186.571
$ws.Range("M136").Value = -13376.6469
$ws.Range("N136").Value = -61286.571
$ws.Range("H141").Value = 57528.75
$ws.Range("J141").Value = 57528.75
$ws.Range("L141").Value = 57528.75
$ws.Range("N141").Value = -67888.75

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 50000844
$ws.Range("I107").Value = 71429230
$ws.Range("J107").Value = 1265
$ws.Range("K107").Value = 214287690
$ws.Range("L107").Value = 3795
$ws.Range("M107").Value = -214285770
$ws.Range("N107").Value = -7635
$ws.Range("H122").Value = 1261.1428
$ws.Range("I122").Value = 955.6667
$ws.Range("J122").Value = 1811
$ws.Range("K122").Value = 2867.0001
$ws.Range("L122").Value = 5433
$ws.Range("M122").Value = -417.0001000000002
$ws.Range("N122").Value = -10333
$ws.Range("H136").Value = 2806.1794
$ws.Range("I136").Value = 3890.647
$ws.Range("J136").Value = 1968.1818
$ws.Range("K136").Value = 11671.941
$ws.Range("L136").Value = 5904.5454
$ws.Range("M136").Value = -9121.940999999999
$ws.Range("N136").Value = -11004.5454
